# ---------------------------------------------------------------------------
# PlayerPerformance_5661.xlsx update
#   1. Insert a new "Player Info" worksheet as the first sheet, with basic
#      player metadata (ID, NAME, BATTING_HAND, BOWL_STYLE).
#   2. In the existing "ODI Batting" sheet, rename the MATCH_CARD_LINK
#      column to MATCH_CODE and replace the full scorecard URLs with just
#      the numeric match code.
#   3. Same rename/value simplification in the "ODI Bowling" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Move($wb.Worksheets.Item(1))
$playerInfo.Name = "Player Info"

# NOTE: sheet references obtained before the Add()/Move() calls above can
# become stale (they resolve by position, and Add() inserts at position 1),
# so look the other two sheets up fresh, *after* the new sheet is in place.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header styling used elsewhere in the
# workbook.
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row - keep everything as text, matching the inline-string cells used
# throughout the rest of the workbook (e.g. the numeric-looking "5661" id).
# Force text storage via a temporary "@" number format, then drop back to
# the default "Normal" style so no stray style index is left on the cell.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5661"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Matthew William Parkinson"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{ "D2" = "4401"; "D3" = "4405"; "D4" = "4472"; "D5" = "4473"; "D6" = "4476" }
foreach ($addr in $battingCodes.Keys) {
    $cell = $battingSheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$addr]
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet - MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{ "B2" = "4401"; "B3" = "4405"; "B4" = "4472"; "B5" = "4473"; "B6" = "4476" }
foreach ($addr in $bowlingCodes.Keys) {
    $cell = $bowlingSheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$addr]
    $cell.Style = "Normal"
}
